$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.133.53'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '1.906.28'
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("E4").Value = '  -0.50%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4615'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.34%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3885'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.54%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07858'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9902'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.98'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.59%  '
$ws.Range("D12").Value = '1.877.94'
$ws.Range("E12").Value = '  -1.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.041'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.749'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07023'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.76%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009926'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.07'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9997'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("D21").Value = '29.136.36'
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.323'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.11'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = '2.117.99'
$ws.Range("E24").Value = '  -0.40%  '
$ws.Range("E25").Value = '  +1.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.99'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.39'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.911'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.56%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '118.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.876'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09317'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8955'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.228'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.320'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.141'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05780'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.169'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02087'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9986'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5702'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.650'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1808'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.725'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.85'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5354'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.000002780'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +63.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.164'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06981'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.839'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.547'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.55%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '113.03'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.37%  '
